$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# JGI "food_bout" app rework: the date/time typed fields become plain text
# fields. "type" is column C.
$ws.Range("C2").Value = "text"   # FB_FOL_date            (was "date")
$ws.Range("C4").Value = "text"   # FB_begin_feed_time      (was "time")
$ws.Range("C5").Value = "text"   # FB_end_feed_time        (was "time")

# Update the active cell/selection on the survey sheet.
$ws.Range("C9").Select()

# Persist so the shared-string table is recompacted (drops the now-unused
# "date"/"time" entries) just like a real Excel save.
$wb.Save()
